# user registration - caso de prueba del nombre (criterios y tamaño)
# Adds a new worksheet "Hoja5" with two more registration test rows
# (one exercising a very long "name" value) and tidies up the
# previously-active sheets' saved selections.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Hoja1: move the saved selection to A16 (single cell, no range).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Activate() | Out-Null
$ws1.Range("A16").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Hoja2: reselect A1:D3 so the active cell resets back to A1.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Activate() | Out-Null
$ws2.Range("A1:D3").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Hoja4: move the saved selection to B2.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Hoja4")
$ws4.Activate() | Out-Null
$ws4.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. New sheet "Hoja5" placed after Hoja4, with two new test rows.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "Hoja5"

$ws5.Range("A1").Value = "Nombre"
$ws5.Range("B1").Value = "Email"
$ws5.Range("C1").Value = "Username"
$ws5.Range("D1").Value = "Contraseña"

$ws5.Range("B2").Value = "validmail011@outlook.com"
$ws5.Range("C2").Value = "jonreyusr010"
$ws5.Range("D2").Value = "P4ssword."

$longName = "nicppdihyvpnokiardkvwlguymkabqzzlafpqszjwbexngljpkfqdbwsddnpkwmlpfsybljiiwcoxbljyogbaifwnvseqlvxvngdjwotcwdgwssxsvngclquzloafjzsrtufcgjjdsngvrknmvnrrvcuwpfoiyfugxkmwrukutqjmajievheoeezmabzlqojexkyyrnadbzkxqsqqltivskigfgiigthbktpcetllhjjqmhxlhxrcuageakznxpbzz"
$ws5.Range("A3").Value = $longName
$ws5.Range("B3").Value = "validmail012@outlook.com"
$ws5.Range("C3").Value = "jonescusr011"
$ws5.Range("D3").Value = "P4ssword."

# Hyperlink the two email cells (mailto:) and restore the shared
# "Hyperlink" cell style afterwards.
$ws5.Hyperlinks.Add($ws5.Range("B2"), "mailto:validmail011@outlook.com") | Out-Null
$ws5.Range("B2").Style = "Hyperlink"

$ws5.Hyperlinks.Add($ws5.Range("B3"), "mailto:validmail012@outlook.com") | Out-Null
$ws5.Range("B3").Style = "Hyperlink"

# Hoja5 becomes the active/selected tab with the cursor on B6.
$ws5.Activate() | Out-Null
$ws5.Range("B6").Select() | Out-Null
